$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30 - this shifts the existing rows 30..47
# down to 31..48 (matching the weekly roll of historical price records).
$ws.Rows.Item(30).Insert()

# Populate the new row 30 with this week's record.
$ws.Range("A30").Value = 10
$ws.Range("B30").Value = "Vega Modelo de Temuco"
$ws.Range("C30").Value = "La Araucanía"
$ws.Range("D30").Value = 44603
$ws.Range("E30").Value = 9
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100107
$ws.Range("H30").Value = "Otros"
$ws.Range("I30").Value = 100107011
$ws.Range("J30").Value = "Tuna"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 40
$ws.Range("N30").Value = 17000
$ws.Range("O30").Value = 18000
$ws.Range("P30").Value = 17500
$ws.Range("Q30").Value = "`$/caja 16 kilos"
$ws.Range("R30").Value = "Provincia de Los Andes"
$ws.Range("S30").Value = 1094
$ws.Range("T30").Value = 16
